$d = $word.ActiveDocument

$replacements = @(
    ,@("2025-07-14 Monday", "2025-07-15 Tuesday")
    ,@("76-53=23", "59-47=12")
    ,@("33+1=34", "58-30=28")
    ,@("22+30=52", "7+19=26")
    ,@("44+30=74", "60-48=12")
    ,@("22-5=17", "76+20=96")
    ,@("69+11=80", "23+10=33")
    ,@("17+19=36", "62-42=20")
    ,@("64+3=67", "76-57=19")
    ,@("17+37=54", "21+30=51")
    ,@("67+10=77", "21+20=41")
    ,@("6+31=37", "49+38=87")
    ,@("2+33=35", "63-60=3")
    ,@("37+27=64", "59+37=96")
    ,@("79+15=94", "90-80=10")
    ,@("2+58=60", "58+7=65")
    ,@("86-15=71", "55+19=74")
    ,@("72+24=96", "47-38=9")
    ,@("10+43=53", "32+33=65")
    ,@("12+12=24", "22+66=88")
    ,@("94-62=32", "68-27=41")
    ,@("30+10=40", "50+26=76")
    ,@("16+62=78", "42-3=39")
    ,@("96-46=50", "1+93=94")
    ,@("51-34=17", "4+27=31")
    ,@("81-52=29", "24+63=87")
    ,@("14+3=17", "40-37=3")
    ,@("74-18=56", "92-29=63")
    ,@("91-57=34", "0+31=31")
    ,@("72+15=87", "10+70=80")
    ,@("29+43=72", "43-22=21")
    ,@("44+27=71", "9-7=2")
    ,@("11+56=67", "34-31=3")
    ,@("67-3=64", "62+31=93")
    ,@("16+64=80", "75+0=75")
    ,@("33+31=64", "38+48=86")
    ,@("78+7=85", "60+36=96")
    ,@("30+22=52", "58-23=35")
    ,@("81+11=92", "54+11=65")
    ,@("68-24=44", "51-18=33")
    ,@("54+3=57", "55-14=41")
    ,@("11-3=8", "31-28=3")
    ,@("71-25=46", "35-9=26")
    ,@("60+10=70", "79-62=17")
    ,@("33+6=39", "13+49=62")
    ,@("93-80=13", "91-9=82")
    ,@("34-4=30", "54+16=70")
    ,@("5+51=56", "28+3=31")
    ,@("4+53=57", "95-50=45")
    ,@("12-5=7", "65-44=21")
    ,@("52+15=67", "58+19=77")
    ,@("88-37=51", "25-25=0")
    ,@("9+82=91", "34+0=34")
    ,@("42+17=59", "17+50=67")
    ,@("64-39=25", "14+20=34")
    ,@("0+32=32", "83-78=5")
    ,@("99-54=45", "20+30=50")
    ,@("66-2=64", "8+83=91")
    ,@("20+33=53", "23+68=91")
    ,@("35+12=47", "96-40=56")
    ,@("28-3=25", "91-89=2")
    ,@("68-36=32", "67-51=16")
    ,@("24+18=42", "37-14=23")
    ,@("57+14=71", "46-0=46")
    ,@("1+29=30", "85-47=38")
    ,@("49-13=36", "67-8=59")
    ,@("60-18=42", "68+16=84")
    ,@("1+31=32", "6+49=55")
    ,@("66-33=33", "38-10=28")
    ,@("8+31=39", "13+20=33")
    ,@("36+42=78", "73+14=87")
    ,@("5+15=20", "67+31=98")
    ,@("16+33=49", "77+2=79")
    ,@("54-31=23", "75-74=1")
    ,@("58-24=34", "73-27=46")
    ,@("78-76=2", "55+15=70")
    ,@("98-79=19", "41+25=66")
    ,@("26+44=70", "2+31=33")
    ,@("76-14=62", "63-11=52")
    ,@("20+61=81", "60-39=21")
    ,@("7+46=53", "61+4=65")
    ,@("26+53=79", "68+14=82")
    ,@("21+16=37", "27+19=46")
    ,@("99-44=55", "50+4=54")
    ,@("58-32=26", "45-25=20")
    ,@("12+25=37", "90-73=17")
    ,@("8+88=96", "86-32=54")
    ,@("71-68=3", "67+24=91")
    ,@("28+60=88", "22+47=69")
    ,@("82-45=37", "37+14=51")
    ,@("12+21=33", "96-6=90")
    ,@("98-52=46", "84-31=53")
    ,@("23+60=83", "75-64=11")
    ,@("35+8=43", "68+2=70")
    ,@("98-54=44", "75-70=5")
    ,@("74-42=32", "35+13=48")
    ,@("46-4=42", "54+31=85")
    ,@("40+30=70", "82-56=26")
    ,@("70-60=10", "8+82=90")
    ,@("58+33=91", "9+19=28")
    ,@("6+5=11", "49-19=30")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Host "Replacements applied: $($replacements.Count)"
